$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047396267"
$ws.Range("D16").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E16").Value = "1904"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 781242
$ws.Range("C17").Value = "73213618"
$ws.Range("D17").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E17").Value = "1904"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 828116
$ws.Range("C18").Value = "1047396267"
$ws.Range("D18").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E18").Value = "1905"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 781242
$ws.Range("C19").Value = "73213618"
$ws.Range("D19").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E19").Value = "1905"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 828116
$ws.Range("C20").Value = "1047396267"
$ws.Range("D20").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E20").Value = "1906"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242
$ws.Range("C21").Value = "73213618"
$ws.Range("D21").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E21").Value = "1906"
$ws.Range("F21").Value = 33125
$ws.Range("G21").Value = 828116
$ws.Range("C22").Value = "1047396267"
$ws.Range("D22").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E22").Value = "1907"
$ws.Range("F22").Value = 31249
$ws.Range("G22").Value = 781242
$ws.Range("C23").Value = "73213618"
$ws.Range("D23").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E23").Value = "1907"
$ws.Range("F23").Value = 33125
$ws.Range("G23").Value = 828116
$ws.Range("C24").Value = "1047396267"
$ws.Range("D24").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E24").Value = "1908"
$ws.Range("F24").Value = 31249
$ws.Range("G24").Value = 781242
$ws.Range("C25").Value = "73213618"
$ws.Range("D25").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E25").Value = "1908"
$ws.Range("F25").Value = 33125
$ws.Range("G25").Value = 828116
$ws.Range("C26").Value = "1047396267"
$ws.Range("D26").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E26").Value = "1909"
$ws.Range("F26").Value = 31249
$ws.Range("G26").Value = 781242
$ws.Range("C27").Value = "73213618"
$ws.Range("D27").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E27").Value = "1909"
$ws.Range("F27").Value = 33125
$ws.Range("G27").Value = 828116
$ws.Range("C28").Value = "1047396267"
$ws.Range("D28").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E28").Value = "1910"
$ws.Range("F28").Value = 31249
$ws.Range("G28").Value = 781242
$ws.Range("C29").Value = "73213618"
$ws.Range("D29").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E29").Value = "1910"
$ws.Range("F29").Value = 33125
$ws.Range("G29").Value = 828116
$ws.Range("C30").Value = "1047396267"
$ws.Range("D30").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E30").Value = "1911"
$ws.Range("F30").Value = 31249
$ws.Range("G30").Value = 781242
$ws.Range("C31").Value = "73213618"
$ws.Range("D31").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E31").Value = "1911"
$ws.Range("F31").Value = 33125
$ws.Range("G31").Value = 828116
$ws.Range("C32").Value = "1047396267"
$ws.Range("D32").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E32").Value = "1912"
$ws.Range("F32").Value = 31249
$ws.Range("G32").Value = 781242
$ws.Range("C33").Value = "73213618"
$ws.Range("D33").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E33").Value = "1912"
$ws.Range("F33").Value = 33125
$ws.Range("G33").Value = 828116
$ws.Range("C34").Value = "1047396267"
$ws.Range("D34").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E34").Value = "2001"
$ws.Range("F34").Value = 31249
$ws.Range("G34").Value = 781242
$ws.Range("C35").Value = "73213618"
$ws.Range("D35").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E35").Value = "2001"
$ws.Range("F35").Value = 33125
$ws.Range("G35").Value = 828116
$ws.Range("C36").Value = "1047396267"
$ws.Range("D36").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E36").Value = "2002"
$ws.Range("F36").Value = 31249
$ws.Range("G36").Value = 781242
$ws.Range("C37").Value = "73213618"
$ws.Range("D37").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E37").Value = "2002"
$ws.Range("F37").Value = 33125
$ws.Range("G37").Value = 828116
$ws.Range("C38").Value = "1047396267"
$ws.Range("D38").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E38").Value = "2003"
$ws.Range("F38").Value = 31249
$ws.Range("G38").Value = 781242
$ws.Range("C39").Value = "73213618"
$ws.Range("D39").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E39").Value = "2003"
$ws.Range("F39").Value = 33125
$ws.Range("G39").Value = 828116
$ws.Range("C40").Value = "1047396267"
$ws.Range("D40").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E40").Value = "2004"
$ws.Range("F40").Value = 31249
$ws.Range("G40").Value = 781242
$ws.Range("C41").Value = "73213618"
$ws.Range("D41").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E41").Value = "2004"
$ws.Range("F41").Value = 33125
$ws.Range("G41").Value = 828116
$ws.Range("C42").Value = "1047396267"
$ws.Range("D42").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E42").Value = "2005"
$ws.Range("F42").Value = 31249
$ws.Range("G42").Value = 781242
$ws.Range("C43").Value = "73213618"
$ws.Range("D43").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E43").Value = "2005"
$ws.Range("F43").Value = 33125
$ws.Range("G43").Value = 828116
$ws.Range("C44").Value = "1047396267"
$ws.Range("D44").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E44").Value = "2006"
$ws.Range("F44").Value = 31249
$ws.Range("G44").Value = 781242
$ws.Range("C45").Value = "73213618"
$ws.Range("D45").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E45").Value = "2006"
$ws.Range("F45").Value = 33125
$ws.Range("G45").Value = 828116
$ws.Range("C46").Value = "1047396267"
$ws.Range("D46").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E46").Value = "2007"
$ws.Range("F46").Value = 31249
$ws.Range("G46").Value = 781242
$ws.Range("C47").Value = "73213618"
$ws.Range("D47").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E47").Value = "2007"
$ws.Range("F47").Value = 33125
$ws.Range("G47").Value = 828116
$ws.Range("C48").Value = "1047396267"
$ws.Range("D48").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E48").Value = "2008"
$ws.Range("F48").Value = 31249
$ws.Range("G48").Value = 781242
$ws.Range("C49").Value = "73213618"
$ws.Range("D49").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E49").Value = "2008"
$ws.Range("F49").Value = 33125
$ws.Range("G49").Value = 828116
$ws.Range("C50").Value = "1047396267"
$ws.Range("D50").Value = "YURIS ZAPATEIRO GUZMAN"
$ws.Range("E50").Value = "2009"
$ws.Range("F50").Value = 29166
$ws.Range("G50").Value = 781242
$ws.Range("C51").Value = "73213618"
$ws.Range("D51").Value = "JESUS MARIA DE LA ROSA PEREZ"
$ws.Range("E51").Value = "2009"
$ws.Range("F51").Value = 30916
$ws.Range("G51").Value = 828116